$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update reporting period (row 8): Q3 2022 -> Q4 2022 ---
# Fecha de inicio del periodo que se informa
$ws.Range("B8").Value = (Get-Date -Year 2022 -Month 10 -Day 1 -Hour 0 -Minute 0 -Second 0)
# Fecha de termino del periodo que se informa
$ws.Range("C8").Value = (Get-Date -Year 2022 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0)
# Fecha de validacion / Fecha de actualizacion
$ws.Range("S8").Value = (Get-Date -Year 2023 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("T8").Value = (Get-Date -Year 2023 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0)

# --- Remove wrap-text formatting on the merged description header (G3:I3) so it
#     matches the plain header style already used by the other title cells
#     (A3/D3 for the fill cell, B3/C3/E3/F3 for the blank continuation cells). ---
$ws.Range("A3").Copy()
$ws.Range("G3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B3").Copy()
$ws.Range("H3:I3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 3 no longer needs the tall custom height once the text isn't wrapping;
# AutoFit drops the explicit override back to the sheet's default height.
$ws.Rows.Item(3).AutoFit()

# --- Column U width adjustment (~55.86 chars; engine quantizes to 1/6-char
#     steps, so 55 is the input that lands closest to the target width) ---
$ws.Columns.Item(21).ColumnWidth = 55
